{"js": "// Replicates the template re-save seen in the diff:\n//  - paragraph 1: the three spell-checked runs (\"A simple \" / \"demonstration\" /\n//    \" of a \" / \"query\" / \" :\") collapse into a single run with the same text\n//    (the <w:proofErr> spell-check markers disappear too).\n//  - paragraph 2 (the \"m:self.name\" field): the scattered <w:instrText> runs\n//    (one of which carried an orange theme color on \"self\") collapse into a\n//    single trimmed instruction \"m:self.name\", and the field gains a\n//    <w:fldChar w:fldCharType=\"separate\"/> marking it as updated/evaluated\n//    (with an empty result, since the referenced login does not exist).\n//  - paragraph 3 (\"End of demonstration.\") and the trailing empty paragraph\n//    keep their text/emptiness unchanged.\n//\n// We rebuild the two affected paragraphs via insertOoxml(..., replace) so the\n// run/field-character structure matches exactly, rather than trying to coax\n// the same result out of higher level text APIs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst OOXML_NS =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>';\nconst OOXML_END =\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n// Paragraph 1: \"A simple demonstration of a query\\u00A0:\" as a single run\n// (note: a non-breaking space, not a plain space, precedes the colon in the\n// original document -- preserved here, only the run split disappears).\nconst introText = \"A simple demonstration of a query\\u00A0:\";\nconst introParagraph =\n  OOXML_NS +\n  '<w:p>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">' + introText + '</w:t></w:r>' +\n  '</w:p>' +\n  OOXML_END;\n\n// Paragraph 2: the \"m:self.name\" field, now updated (begin / instrText /\n// separate / end) instead of left purely as instruction text.\nconst fieldParagraph =\n  OOXML_NS +\n  '<w:p>' +\n  '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n  '<w:r><w:instrText xml:space=\"preserve\">m:self.name</w:instrText></w:r>' +\n  '<w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n  '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n  '</w:p>' +\n  OOXML_END;\n\nparagraphs.items[0].insertOoxml(introParagraph, Word.InsertLocation.replace);\nparagraphs.items[1].insertOoxml(fieldParagraph, Word.InsertLocation.replace);\nawait context.sync();\n\n// The re-saved sectPr also reports header/footer distance 0 instead of the\n// original 708 twips (35.4pt) -- harmless since the document has no header\n// or footer content, but it's part of the same re-save; mirror it too.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\nconst pageSetup = sections.items[0].pageSetup;\npageSetup.headerDistance = 0;\npageSetup.footerDistance = 0;\nawait context.sync();\n", "ps1": "# Replicates the template re-save seen in the diff:\n#  - paragraph 1: the spell-checked runs (\"A simple \" / \"demonstration\" /\n#    \" of a \" / \"query\" / \" :\") collapse into a single run with the same\n#    text (the spell-check <w:proofErr> markers disappear too).\n#  - paragraph 2 (the \"m:self.name\" field): the scattered <w:instrText> runs\n#    (one of which carried an orange theme color on \"self\") collapse into a\n#    single trimmed instruction \"m:self.name\", and the field gains a\n#    <w:fldChar w:fldCharType=\"separate\"/> marking it as updated/evaluated\n#    (with an empty result, since the referenced login does not exist).\n#  - paragraph 3 (\"End of demonstration.\") and the trailing empty paragraph\n#    keep their text/emptiness unchanged.\n#\n# We rebuild the two affected paragraphs with Range.InsertXML(...) (which\n# replaces the range content with the supplied WordprocessingML) so the\n# run/field-character structure matches exactly, rather than trying to coax\n# the same result out of higher level text APIs.\n\n$d = $word.ActiveDocument\n\n$ooxmlNs = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$ooxmlEnd = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Paragraph 1: \"A simple demonstration of a query<NBSP>:\" as a single run\n# (note: a non-breaking space, not a plain space, precedes the colon in the\n# original document -- preserved here, only the run split disappears).\n$introText = \"A simple demonstration of a query$([char]0x00A0):\"\n$introParagraph = $ooxmlNs +\n  '<w:p>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">' + $introText + '</w:t></w:r>' +\n  '</w:p>' +\n  $ooxmlEnd\n\n# Paragraph 2: the \"m:self.name\" field, now updated (begin / instrText /\n# separate / end) instead of left purely as instruction text.\n$fieldParagraph = $ooxmlNs +\n  '<w:p>' +\n  '<w:r><w:fldChar w:fldCharType=\"begin\"/></w:r>' +\n  '<w:r><w:instrText xml:space=\"preserve\">m:self.name</w:instrText></w:r>' +\n  '<w:r><w:fldChar w:fldCharType=\"separate\"/></w:r>' +\n  '<w:r><w:fldChar w:fldCharType=\"end\"/></w:r>' +\n  '</w:p>' +\n  $ooxmlEnd\n\n$p1 = $d.Paragraphs.Item(1).Range\n[void]$p1.InsertXML($introParagraph)\n\n$p2 = $d.Paragraphs.Item(2).Range\n[void]$p2.InsertXML($fieldParagraph)\n\n# The re-saved sectPr also reports header/footer distance 0 instead of the\n# original 708 twips (35.4pt) -- harmless since the document has no header\n# or footer content, but it's part of the same re-save; mirror it too.\n$ps = $d.PageSetup\n$ps.HeaderDistance = 0\n$ps.FooterDistance = 0\n"}
